$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, shifting existing rows 91.. down to 92..
$ws.Rows("91:91").Insert()

# Populate the newly inserted row 91 with the new data record
$ws.Cells.Item(91, 1).Value = 4
$ws.Cells.Item(91, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(91, 3).Value = "Los Lagos"
$ws.Cells.Item(91, 4).Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(44651)
$ws.Cells.Item(91, 5).Value = 10
$ws.Cells.Item(91, 6).Value = 100112009
$ws.Cells.Item(91, 7).Value = "Acelga"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 20
$ws.Cells.Item(91, 11).Value = 10000
$ws.Cells.Item(91, 12).Value = 10000
$ws.Cells.Item(91, 13).Value = 10000
$ws.Cells.Item(91, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(91, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(91, 16).Value = 833
$ws.Cells.Item(91, 17).Value = 12
$ws.Cells.Item(91, 18).Value = "Hortaliza"
